$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column C (year) as Text so numeric-looking years stay text,
# matching the inlineStr type used throughout the rest of the column.
$ws.Range("C21:C56").NumberFormat = "@"

$ws.Cells.Item(21, 1).Value = "Promotie"
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(21, 3).Value = "1963"
$ws.Cells.Item(21, 4).Value = 20

$ws.Cells.Item(22, 1).Value = "Promotie"
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = "1977"
$ws.Cells.Item(22, 4).Value = 20

$ws.Cells.Item(23, 1).Value = "Promotie"
$ws.Cells.Item(23, 2).Value = 2
$ws.Cells.Item(23, 3).Value = "1978"
$ws.Cells.Item(23, 4).Value = 20

$ws.Cells.Item(24, 1).Value = "Promotie"
$ws.Cells.Item(24, 2).Value = 5
$ws.Cells.Item(24, 3).Value = "1980"
$ws.Cells.Item(24, 4).Value = 20

$ws.Cells.Item(25, 1).Value = "Promotie"
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = "1981"
$ws.Cells.Item(25, 4).Value = 20

$ws.Cells.Item(26, 1).Value = "Promotie"
$ws.Cells.Item(26, 2).Value = 5
$ws.Cells.Item(26, 3).Value = "1982"
$ws.Cells.Item(26, 4).Value = 20

$ws.Cells.Item(27, 1).Value = "Promotie"
$ws.Cells.Item(27, 2).Value = 4
$ws.Cells.Item(27, 3).Value = "1983"
$ws.Cells.Item(27, 4).Value = 20

$ws.Cells.Item(28, 1).Value = "Promotie"
$ws.Cells.Item(28, 2).Value = 10
$ws.Cells.Item(28, 3).Value = "1984"
$ws.Cells.Item(28, 4).Value = 20

$ws.Cells.Item(29, 1).Value = "Promotie"
$ws.Cells.Item(29, 2).Value = 13
$ws.Cells.Item(29, 3).Value = "1985"
$ws.Cells.Item(29, 4).Value = 20

$ws.Cells.Item(30, 1).Value = "Promotie"
$ws.Cells.Item(30, 2).Value = 8
$ws.Cells.Item(30, 3).Value = "1986"
$ws.Cells.Item(30, 4).Value = 20

$ws.Cells.Item(31, 1).Value = "Promotie"
$ws.Cells.Item(31, 2).Value = 16
$ws.Cells.Item(31, 3).Value = "1987"
$ws.Cells.Item(31, 4).Value = 20

$ws.Cells.Item(32, 1).Value = "Promotie"
$ws.Cells.Item(32, 2).Value = 17
$ws.Cells.Item(32, 3).Value = "1988"
$ws.Cells.Item(32, 4).Value = 20

$ws.Cells.Item(33, 1).Value = "Promotie"
$ws.Cells.Item(33, 2).Value = 20
$ws.Cells.Item(33, 3).Value = "1989"
$ws.Cells.Item(33, 4).Value = 20

$ws.Cells.Item(34, 1).Value = "Promotie"
$ws.Cells.Item(34, 2).Value = 11
$ws.Cells.Item(34, 3).Value = "1990"
$ws.Cells.Item(34, 4).Value = 20

$ws.Cells.Item(35, 1).Value = "Promotie"
$ws.Cells.Item(35, 2).Value = 15
$ws.Cells.Item(35, 3).Value = "1991"
$ws.Cells.Item(35, 4).Value = 20

$ws.Cells.Item(36, 1).Value = "Promotie"
$ws.Cells.Item(36, 2).Value = 20
$ws.Cells.Item(36, 3).Value = "1992"
$ws.Cells.Item(36, 4).Value = 20

$ws.Cells.Item(37, 1).Value = "Promotie"
$ws.Cells.Item(37, 2).Value = 11
$ws.Cells.Item(37, 3).Value = "1993"
$ws.Cells.Item(37, 4).Value = 20

$ws.Cells.Item(38, 1).Value = "Promotie"
$ws.Cells.Item(38, 2).Value = 17
$ws.Cells.Item(38, 3).Value = "1994"
$ws.Cells.Item(38, 4).Value = 20

$ws.Cells.Item(39, 1).Value = "Promotoe"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = "1994"
$ws.Cells.Item(39, 4).Value = 20

$ws.Cells.Item(40, 1).Value = "Promotie"
$ws.Cells.Item(40, 2).Value = 9
$ws.Cells.Item(40, 3).Value = "1995"
$ws.Cells.Item(40, 4).Value = 20

$ws.Cells.Item(41, 1).Value = "Promotie"
$ws.Cells.Item(41, 2).Value = 7
$ws.Cells.Item(41, 3).Value = "1996"
$ws.Cells.Item(41, 4).Value = 20

$ws.Cells.Item(42, 1).Value = "Promotie"
$ws.Cells.Item(42, 2).Value = 11
$ws.Cells.Item(42, 3).Value = "1997"
$ws.Cells.Item(42, 4).Value = 20

$ws.Cells.Item(43, 1).Value = "Promotie"
$ws.Cells.Item(43, 2).Value = 5
$ws.Cells.Item(43, 3).Value = "1998"
$ws.Cells.Item(43, 4).Value = 20

$ws.Cells.Item(44, 1).Value = "Promotie"
$ws.Cells.Item(44, 2).Value = 4
$ws.Cells.Item(44, 3).Value = "1999"
$ws.Cells.Item(44, 4).Value = 20

$ws.Cells.Item(45, 1).Value = "Promotie"
$ws.Cells.Item(45, 2).Value = 9
$ws.Cells.Item(45, 3).Value = "2000"
$ws.Cells.Item(45, 4).Value = 21

$ws.Cells.Item(46, 1).Value = "Promotie"
$ws.Cells.Item(46, 2).Value = 1
$ws.Cells.Item(46, 3).Value = "2001"
$ws.Cells.Item(46, 4).Value = 21

$ws.Cells.Item(47, 1).Value = "Promotie"
$ws.Cells.Item(47, 2).Value = 4
$ws.Cells.Item(47, 3).Value = "2002"
$ws.Cells.Item(47, 4).Value = 21

$ws.Cells.Item(48, 1).Value = "Promotie"
$ws.Cells.Item(48, 2).Value = 5
$ws.Cells.Item(48, 3).Value = "2003"
$ws.Cells.Item(48, 4).Value = 21

$ws.Cells.Item(49, 1).Value = "Promotie"
$ws.Cells.Item(49, 2).Value = 3
$ws.Cells.Item(49, 3).Value = "2004"
$ws.Cells.Item(49, 4).Value = 21

$ws.Cells.Item(50, 1).Value = "Promotie"
$ws.Cells.Item(50, 2).Value = 1
$ws.Cells.Item(50, 3).Value = "2005"
$ws.Cells.Item(50, 4).Value = 21

$ws.Cells.Item(51, 1).Value = "Promotie"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = "2006"
$ws.Cells.Item(51, 4).Value = 21

$ws.Cells.Item(52, 1).Value = "Promotie"
$ws.Cells.Item(52, 2).Value = 3
$ws.Cells.Item(52, 3).Value = "2007"
$ws.Cells.Item(52, 4).Value = 21

$ws.Cells.Item(53, 1).Value = "Promotie"
$ws.Cells.Item(53, 2).Value = 3
$ws.Cells.Item(53, 3).Value = "2008"
$ws.Cells.Item(53, 4).Value = 21

$ws.Cells.Item(54, 1).Value = "Promotie"
$ws.Cells.Item(54, 2).Value = 2
$ws.Cells.Item(54, 3).Value = "2010"
$ws.Cells.Item(54, 4).Value = 21

$ws.Cells.Item(55, 1).Value = "Promotie"
$ws.Cells.Item(55, 2).Value = 2
$ws.Cells.Item(55, 3).Value = "2011"
$ws.Cells.Item(55, 4).Value = 21

$ws.Cells.Item(56, 1).Value = "Promotie"
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(56, 3).Value = "2012"
$ws.Cells.Item(56, 4).Value = 21
